# Rebuild the "execution graph" draft on Sheet1.
#
# The previous layout (a chain of helper formulas spread across A1:G9)
# is replaced by a smaller graph:
#   - A1 = 5 - IF(I1=I2, I3, I4)            (decision driven by I1:I4)
#   - B1..B6 hold a separate small chain culminating in D9 = COUNT(D1:D8)
#   - I1..I4 hold the lookup values/formula feeding A1's IF()
#
# Clear the old used range first so stale cells (C1, D1:G1, C2:G2, D3,
# E3, F4, G4, F5, A6, A7, B7, C7, B8, B9, ...) are removed rather than
# merely overwritten, then write the new cells/formulas one by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:I9").ClearContents()

# Row 1
$ws.Range("A1").Formula = "=5-IF(I1=I2,I3,I4)"
$ws.Range("B1").Value = 10
$ws.Range("I1").Value = 1

# Row 2
$ws.Range("B2").Value = 5
$ws.Range("I2").Value = 2

# Row 3
$ws.Range("B3").Value = 8
$ws.Range("I3").Value = 3

# Row 4
$ws.Range("B4").Formula = "=D9"
$ws.Range("I4").Formula = "=B5+B5"

# Row 5
$ws.Range("B5").Formula = "=D5-E5"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 3

# Row 6
$ws.Range("B6").Value = 0

# Row 9
$ws.Range("D9").Formula = "=COUNT(D1:D8)"

# Match the saved cursor/selection position.
$null = $ws.Range("D18").Select()
